# Apply the EPEX Spot workbook update:
#  - "Prix Spot": insert a new date column ("03-nov") before the existing
#    "01-oct." column (column DH), shifting DH:EL -> DI:EM, filling the new
#    column with "-" placeholders for the data rows.
#  - "Gaz" and "CO2": append two new daily rows (2025-11-01 / 2025-11-02).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Prix Spot" — insert a new column at DH (column 112).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Shifts the existing DH:EL (01-oct. .. 31-oct.) block one column to the
# right (-> DI:EM), leaving a blank column DH ready to receive the new date.
$ws1.Columns.Item(112).Insert()

# Header for the newly inserted column.
$ws1.Cells.Item(1,112).Value = "03-nov"

# The data rows (2..25) for that new day have no recorded price yet.
$ws1.Range($ws1.Cells.Item(2,112), $ws1.Cells.Item(25,112)).Value = "-"

# ---------------------------------------------------------------------------
# Sheet 2: "Gaz" — append the two following calendar days.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Force column A to stay text (otherwise Excel auto-converts the
# "YYYY-MM-DD" strings into date values), matching the existing rows.
$rngA2 = $ws2.Range("A140:A141")
$rngA2.NumberFormat = "@"
$ws2.Cells.Item(140,1).Value = "2025-11-01"
$ws2.Cells.Item(140,2).Value = 29.8
$ws2.Cells.Item(141,1).Value = "2025-11-02"
$ws2.Cells.Item(141,2).Value = 29.8
$rngA2.Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet 3: "CO2" — append the same two calendar days.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$rngA3 = $ws3.Range("A140:A141")
$rngA3.NumberFormat = "@"
$ws3.Cells.Item(140,1).Value = "2025-11-01"
$ws3.Cells.Item(140,2).Value = 78
$ws3.Cells.Item(141,1).Value = "2025-11-02"
$ws3.Cells.Item(141,2).Value = 78
$rngA3.Style = "Normal"
